# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates to the leve profit tables across all class sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 10599.6
$ws.Range("I69").Value = 3999
$ws.Range("J69").Value = 12249.75
$ws.Range("K69").Value = 11997
$ws.Range("L69").Value = 36749.25
$ws.Range("M69").Value = -11123
$ws.Range("N69").Value = -38497.25
$ws.Range("H72").Value = 10599.6
$ws.Range("I72").Value = 3999
$ws.Range("J72").Value = 12249.75
$ws.Range("K72").Value = 35991
$ws.Range("L72").Value = 110247.75
$ws.Range("M72").Value = -31623
$ws.Range("N72").Value = -118983.75
$ws.Range("H132").Value = 3557
$ws.Range("I132").Value = 2462
$ws.Range("K132").Value = 7386
$ws.Range("M132").Value = -4856
$ws.Range("H137").Value = 1703
$ws.Range("I137").Value = 1618.75
$ws.Range("J137").Value = 1787.25
$ws.Range("K137").Value = 4856.25
$ws.Range("L137").Value = 5361.75
$ws.Range("M137").Value = -2306.25
$ws.Range("N137").Value = -10461.75
$ws.Range("H138").Value = 7578533
$ws.Range("I138").Value = 1343.0555
$ws.Range("J138").Value = 10419979
$ws.Range("K138").Value = 4029.1665
$ws.Range("L138").Value = 31259937
$ws.Range("M138").Value = 1110.8335
$ws.Range("N138").Value = -31270217

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3978.758
$ws.Range("I61").Value = 2773.25
$ws.Range("K61").Value = 2773.25
$ws.Range("M61").Value = -2561.25
$ws.Range("H88").Value = 1993.1666
$ws.Range("I88").Value = 1388.8334
$ws.Range("J88").Value = 2597.5
$ws.Range("K88").Value = 1388.8334
$ws.Range("L88").Value = 2597.5
$ws.Range("M88").Value = -982.8334
$ws.Range("N88").Value = -3409.5
$ws.Range("H91").Value = 1993.1666
$ws.Range("I91").Value = 1388.8334
$ws.Range("J91").Value = 2597.5
$ws.Range("K91").Value = 1388.8334
$ws.Range("L91").Value = 2597.5
$ws.Range("M91").Value = 15.16660000000002
$ws.Range("N91").Value = -5405.5
$ws.Range("H132").Value = 3396.889
$ws.Range("I132").Value = 2957.923
$ws.Range("K132").Value = 8873.769
$ws.Range("M132").Value = -6343.769
$ws.Range("H136").Value = 3978.758
$ws.Range("I136").Value = 2773.25
$ws.Range("K136").Value = 8319.75
$ws.Range("M136").Value = -5769.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 9500
$ws.Range("J19").Value = 9500
$ws.Range("L19").Value = 9500
$ws.Range("N19").Value = -9846
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("N109").ClearContents()
$ws.Range("H134").Value = 1879.8906
$ws.Range("I134").Value = 1823.55
$ws.Range("J134").Value = 2725
$ws.Range("K134").Value = 5470.65
$ws.Range("L134").Value = 8175
$ws.Range("M134").Value = -2935.65
$ws.Range("N134").Value = -13245
$ws.Range("H135").Value = 59600
$ws.Range("J135").Value = 59600
$ws.Range("L135").Value = 59600
$ws.Range("N135").Value = -69740

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 23826.777
$ws.Range("I134").Value = 11906.272
$ws.Range("K134").Value = 35718.81600000001
$ws.Range("M134").Value = -33183.81600000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1823.4706
$ws.Range("I114").Value = 1308.875
$ws.Range("J114").Value = 2280.889
$ws.Range("K114").Value = 3926.625
$ws.Range("L114").Value = 6842.667
$ws.Range("M114").Value = -672.625
$ws.Range("N114").Value = -13350.667
$ws.Range("H121").Value = 1764.7778
$ws.Range("I121").Value = 2462.8
$ws.Range("J121").Value = 892.25
$ws.Range("K121").Value = 7388.400000000001
$ws.Range("L121").Value = 2676.75
$ws.Range("M121").Value = -6078.400000000001
$ws.Range("N121").Value = -5296.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2945.2307
$ws.Range("I80").Value = 2598.889
$ws.Range("J80").Value = 3724.5
$ws.Range("K80").Value = 2598.889
$ws.Range("L80").Value = 3724.5
$ws.Range("M80").Value = -1600.889
$ws.Range("N80").Value = -5720.5
$ws.Range("H83").Value = 2945.2307
$ws.Range("I83").Value = 2598.889
$ws.Range("J83").Value = 3724.5
$ws.Range("K83").Value = 12994.445
$ws.Range("L83").Value = 18622.5
$ws.Range("M83").Value = -8002.445
$ws.Range("N83").Value = -28606.5
$ws.Range("H122").Value = 3047.5
$ws.Range("I122").Value = 2427.1
$ws.Range("J122").Value = 6149.5
$ws.Range("K122").Value = 7281.299999999999
$ws.Range("L122").Value = 18448.5
$ws.Range("M122").Value = -4831.299999999999
$ws.Range("N122").Value = -23348.5
$ws.Range("H123").Value = 32481.25
$ws.Range("J123").Value = 32481.25
$ws.Range("L123").Value = 32481.25
$ws.Range("N123").Value = -37381.25
$ws.Range("H132").Value = 4634.1113
$ws.Range("I132").Value = 4285.4165
$ws.Range("K132").Value = 12856.2495
$ws.Range("M132").Value = -10326.2495

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 130000
$ws.Range("I3").Value = 130000
$ws.Range("K3").Value = 130000
$ws.Range("M3").Value = -129888
$ws.Range("H14").Value = 9501
$ws.Range("I14").Value = 8004
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 8004
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -7832
$ws.Range("N14").Value = -10344
$ws.Range("H15").Value = 130000
$ws.Range("I15").Value = 130000
$ws.Range("K15").Value = 130000
$ws.Range("M15").Value = -129830
$ws.Range("H16").Value = 1865.9412
$ws.Range("I16").Value = 1865.9412
$ws.Range("K16").Value = 1865.9412
$ws.Range("M16").Value = -1695.9412
$ws.Range("H22").Value = 2088.647
$ws.Range("I22").Value = 2685.8333
$ws.Range("J22").Value = 1762.909
$ws.Range("K22").Value = 2685.8333
$ws.Range("L22").Value = 1762.909
$ws.Range("M22").Value = -2390.8333
$ws.Range("N22").Value = -2352.909
$ws.Range("H27").Value = 2088.647
$ws.Range("I27").Value = 2685.8333
$ws.Range("J27").Value = 1762.909
$ws.Range("K27").Value = 2685.8333
$ws.Range("L27").Value = 1762.909
$ws.Range("M27").Value = -2578.8333
$ws.Range("N27").Value = -1976.909
$ws.Range("H82").Value = 3455.6667
$ws.Range("I82").Value = 3783
$ws.Range("J82").Value = 2801
$ws.Range("K82").Value = 3783
$ws.Range("L82").Value = 2801
$ws.Range("M82").Value = -3422
$ws.Range("N82").Value = -3523
$ws.Range("H85").Value = 3455.6667
$ws.Range("I85").Value = 3783
$ws.Range("J85").Value = 2801
$ws.Range("K85").Value = 3783
$ws.Range("L85").Value = 2801
$ws.Range("M85").Value = -2535
$ws.Range("N85").Value = -5297
$ws.Range("H108").Value = 33633.332
$ws.Range("J108").Value = 33633.332
$ws.Range("L108").Value = 33633.332
$ws.Range("N108").Value = -41313.332
$ws.Range("H132").Value = 3894.087
$ws.Range("I132").Value = 3034
$ws.Range("J132").Value = 6331
$ws.Range("K132").Value = 9102
$ws.Range("L132").Value = 18993
$ws.Range("M132").Value = -6572
$ws.Range("N132").Value = -24053
$ws.Range("H136").Value = 3628.0688
$ws.Range("I136").Value = 3161.625
$ws.Range("K136").Value = 9484.875
$ws.Range("M136").Value = -6934.875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 250
$ws.Range("I19").Value = 250
$ws.Range("K19").Value = 250
$ws.Range("M19").Value = -76
$ws.Range("H113").Value = 686.55884
$ws.Range("I113").Value = 796.0952
$ws.Range("K113").Value = 2388.2856
$ws.Range("M113").Value = -218.2856000000002

